# Append a new case at the top of the active "案件情報" (job listing) table:
#   - a new row is inserted right before the old row 7 so the sheet keeps its
#     existing "newest first" ordering for this entry's priority score
#   - every data row's "取得日時" (fetched-at) timestamp is refreshed
#   - a fresh set of hyperlinks is rebuilt for column F so the F-column link
#     targets line up 1:1 with the (now shifted) URL text in each row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-15 01:25:52"

# Shift rows 7:17 down to 8:18, duplicating row 7's formatting (incl. the
# hyperlink cell style) into the freshly inserted row.
$ws.Rows("7:7").Insert()

# Fill in the brand-new row 7 with the newly scraped listing.
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "【長期案件】生成AIを利用したチャットボット作成のPMOを募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5471035"
$ws.Range("G7").Value = 310
$ws.Range("H7").Value = "🔥AI,Ai"

# Refresh the fetched-at timestamp on every other data row (2-6 stay in
# place, 8-18 are the previously-existing rows that the insert pushed down).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp
$ws.Range("A13").Value = $newTimestamp
$ws.Range("A14").Value = $newTimestamp
$ws.Range("A15").Value = $newTimestamp
$ws.Range("A16").Value = $newTimestamp
$ws.Range("A17").Value = $newTimestamp
$ws.Range("A18").Value = $newTimestamp

# The row insert already shifted the old hyperlink objects down along with
# their cells, but there is no clean single-cell way to splice a brand-new
# hyperlink into the middle of that collection without leaving a stale
# duplicate behind, so rebuild the whole F2:F18 hyperlink set from scratch
# in row order (matches the final URL text in each cell 1:1).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5450864")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5460294")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5460267")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5471108")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5470737")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5471035")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5471032")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5470814")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5471068")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5470623")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5470812")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5470403")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5471022")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5470150")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5450323")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5470726")
